$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 2836.1927670000005
$ws.Range("M2").Value = 4600.0
$ws.Range("N2").Value = 1780.0

# Row 3
$ws.Range("A3").Value = 7119.401966000001
$ws.Range("M3").Value = 4600.0
$ws.Range("N3").Value = 1780.0

# Row 4
$ws.Range("A4").Value = 5315.985534000003
$ws.Range("M4").Value = 4600.0
$ws.Range("N4").Value = 1780.0

# Row 5
$ws.Range("A5").Value = 2566.5572329999995
$ws.Range("M5").Value = 4600.0
$ws.Range("N5").Value = 1780.0

# Row 6
$ws.Range("A6").Value = 312.6072329999988
$ws.Range("M6").Value = 4600.0
$ws.Range("N6").Value = 1780.0

# Row 7
$ws.Range("M7").Value = 3225.771698999997
$ws.Range("N7").Value = 1780.0

# Row 8
$ws.Range("M8").Value = 1411.614465999999
$ws.Range("N8").Value = 1780.0

# Row 9
$ws.Range("K9").Value = 36000.803301
$ws.Range("M9").Value = 961.1966990000001
$ws.Range("N9").Value = 1780.0

# Row 10
$ws.Range("M10").Value = 577.6269660000034
$ws.Range("N10").Value = 1780.0

# Row 11
$ws.Range("L11").Value = 0.0
$ws.Range("M11").Value = 4600.0
$ws.Range("N11").Value = 1050.990801
